# cryptos list refresh — Mon Apr 17 10:16:25 UTC 2023 (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 47/48 swapped places (BabyDogeCoin now ranks above PancakeSwap) ---
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000361"
$ws.Range("E47").Value = "  -1.87%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.640"
$ws.Range("E48").Value = "  -1.64%  "

# --- Refreshed price / 1h-volume figures for every other row ---
$ws.Range("D2").Value = "'30.086.44"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "'2.104.32"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").Value = "'344.88"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'0.5168"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D8").Value = "'0.4423"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").Value = "'0.09484"
$ws.Range("E9").Value = "  +4.49%  "
$ws.Range("D10").Value = "'52.43"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").Value = "'2.109.49"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "'6.717"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "'8.100"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "'99.64"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "'0.00001167"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "'1.009"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'20.71"
$ws.Range("E19").Value = "  +6.47%  "
$ws.Range("D20").Value = "'0.06710"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "'6.211"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("D23").Value = "'30.162.81"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "'12.71"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "'2.333"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'2.355.91"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'22.06"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").Value = "'164.78"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'2.547"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "'133.71"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").Value = "'1.162"
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").Value = "'0.1056"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").Value = "'1.633"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'6.254"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").Value = "'3.968"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'6.204"
$ws.Range("E36").Value = "  +4.28%  "
$ws.Range("D37").Value = "'10.11"
$ws.Range("E37").Value = "  -3.54%  "
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("D39").Value = "'0.06786"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "'0.2282"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("D41").Value = "'0.6950"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'12.52"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").Value = "'1.311"
$ws.Range("E43").Value = "  +3.84%  "
$ws.Range("D44").Value = "'0.6677"
$ws.Range("E44").Value = "  +3.54%  "
$ws.Range("D45").Value = "'14.17"
$ws.Range("E45").Value = "  -6.11%  "
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D49").Value = "'1.224"
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").Value = "'82.41"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "'0.07198"
$ws.Range("E51").Value = "  -1.73%  "
